$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new vocabulary words (rows 137-205) to column A, mirroring the
# existing List1-List4 layout. Plain entries inherit the sheet default font
# style; the new "List5" section header (row 182) gets the same bold red
# style as the other list headers (copied below from A1).
$ws.Range("A137").Value = "amicable"
$ws.Range("A138").Value = "amity"
$ws.Range("A139").Value = "enmity"
$ws.Range("A140").Value = "enamor"
$ws.Range("A141").Value = "pamper"
$ws.Range("A142").Value = "amphibian"
$ws.Range("A143").Value = "ambition"
$ws.Range("A144").Value = "ambiguous"
$ws.Range("A145").Value = "ambivalent"
$ws.Range("A146").Value = "ambidextrous"
$ws.Range("A147").Value = "ambience"
$ws.Range("A148").Value = "atmosphere"
$ws.Range("A149").Value = "acquiesce"
$ws.Range("A150").Value = "efflorescence"
$ws.Range("A151").Value = "quiescent"
$ws.Range("A152").Value = "quietus"
$ws.Range("A153").Value = "acquittal"
$ws.Range("A154").Value = "convict"
$ws.Range("A155").Value = "addle"
$ws.Range("A156").Value = "baffle"
$ws.Range("A157").Value = "bewilder"
$ws.Range("A158").Value = "adept"
$ws.Range("A159").Value = "adopt"
$ws.Range("A160").Value = "adapt"
$ws.Range("A161").Value = "inept"
$ws.Range("A162").Value = "inapt"
$ws.Range("A163").Value = "deft"
$ws.Range("A164").Value = "adroit"
$ws.Range("A165").Value = "daft"
$ws.Range("A166").Value = "adjacent"
$ws.Range("A167").Value = "reject"
$ws.Range("A168").Value = "deject"
$ws.Range("A169").Value = "abject"
$ws.Range("A170").Value = "conjuecture"
$ws.Range("A171").Value = "ejaculate"
$ws.Range("A172").Value = "abut"
$ws.Range("A173").Value = "adjoin"
$ws.Range("A174").Value = "adjourn"
$ws.Range("A175").Value = "journey"
$ws.Range("A176").Value = "sojourn"
$ws.Range("A177").Value = "diurnal"
$ws.Range("A178").Value = "nocturnal"
$ws.Range("A179").Value = "adulterate"
$ws.Range("A180").Value = "adultery"
$ws.Range("A181").Value = "adulate"
$ws.Range("A182").Value = "List5"
$ws.Range("A183").Value = "annihilate"
$ws.Range("A184").Value = "nihilism"
$ws.Range("A185").Value = "nullify"
$ws.Range("A186").Value = "annul"
$ws.Range("A187").Value = "annual"
$ws.Range("A188").Value = "eliminate"
$ws.Range("A189").Value = "annotate"
$ws.Range("A190").Value = "notation"
$ws.Range("A191").Value = "connotation"
$ws.Range("A192").Value = "denotation"
$ws.Range("A193").Value = "notorious"
$ws.Range("A194").Value = "antagonize"
$ws.Range("A195").Value = "antagonist"
$ws.Range("A196").Value = "protagonist"
$ws.Range("A197").Value = "opponent"
$ws.Range("A198").Value = "rival"
$ws.Range("A199").Value = "agitate"
$ws.Range("A200").Value = "cogitate"
$ws.Range("A201").Value = "cogent"
$ws.Range("A202").Value = "exigent"
$ws.Range("A203").Value = "aquatic"
$ws.Range("A204").Value = "aquarium"
$ws.Range("A205").Value = "Aquarius"

# Match the "ListN" header style (bold red Times New Roman) by copying the
# format from the existing List1 header cell (A1).
$excel.CutCopyMode = $false
$ws.Range("A1").Copy()
$ws.Range("A182").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The other three list headers (List2/List3/List4) pick up an explicit
# 14.25pt row height as a side effect of the edit, matching List1 row.
$ws.Rows(39).RowHeight = 14.25
$ws.Rows(97).RowHeight = 14.25
$ws.Rows(133).RowHeight = 14.25

# Leave the selection where Excel lands after typing the last entry.
$ws.Range("A206").Select()
